$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells that receive numeric-looking text to stay as text
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '25.749.33'
$ws.Range("E2").Value = '  -2.39%  '

$ws.Range("D3").Value = '1.752.79'
$ws.Range("E3").Value = '  -4.27%  '

$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = '236.72'
$ws.Range("E5").Value = '  -5.14%  '

$ws.Range("E6").Value = '  -0.05%  '

$ws.Range("D7").Value = '0.5070'
$ws.Range("E7").Value = '  -3.34%  '

$ws.Range("D8").Value = '41.42'
$ws.Range("E8").Value = '  -6.39%  '

$ws.Range("D9").Value = '0.2646'
$ws.Range("E9").Value = '  -4.22%  '

$ws.Range("D10").Value = '0.06172'
$ws.Range("E10").Value = '  -9.16%  '

$ws.Range("D11").Value = '1.752.30'
$ws.Range("E11").Value = '  -4.37%  '

$ws.Range("D12").Value = '0.06920'
$ws.Range("E12").Value = '  -2.45%  '

$ws.Range("D13").Value = '15.58'
$ws.Range("E13").Value = '  -4.92%  '

$ws.Range("D14").Value = '0.6006'
$ws.Range("E14").Value = '  -12.20%  '

$ws.Range("D15").Value = '4.489'
$ws.Range("E15").Value = '  -7.10%  '

$ws.Range("E16").Value = '  -9.82%  '

$ws.Range("E17").Value = '  -0.06%  '

$ws.Range("E18").Value = '  -0.10%  '

$ws.Range("D19").Value = '25.770.72'
$ws.Range("E19").Value = '  -2.38%  '

$ws.Range("D20").Value = '0.000006850'
$ws.Range("E20").Value = '  -5.90%  '

$ws.Range("D21").Value = '11.68'
$ws.Range("E21").Value = '  -11.06%  '

$ws.Range("D22").Value = '1.978.01'
$ws.Range("E22").Value = '  -4.80%  '

$ws.Range("E23").Value = '  -8.80%  '

$ws.Range("D24").Value = '8.250'
$ws.Range("E24").Value = '  -7.55%  '

$ws.Range("D25").Value = '5.195'
$ws.Range("E25").Value = '  -10.07%  '

$ws.Range("D26").Value = '137.27'
$ws.Range("E26").Value = '  -3.60%  '

$ws.Range("D27").Value = '1.469'
$ws.Range("E27").Value = '  -11.78%  '

$ws.Range("D28").Value = '1.819'
$ws.Range("E28").Value = '  -9.51%  '

$ws.Range("D29").Value = '14.97'
$ws.Range("E29").Value = '  -9.33%  '

$ws.Range("D30").Value = '102.53'
$ws.Range("E30").Value = '  -5.62%  '

$ws.Range("D31").Value = '0.08192'
$ws.Range("E31").Value = '  -5.80%  '

$ws.Range("E32").Value = '  -9.21%  '

$ws.Range("E33").Value = '  -10.03%  '

$ws.Range("D34").Value = '0.04497'
$ws.Range("E34").Value = '  -3.54%  '

$ws.Range("E35").Value = '  -0.09%  '

$ws.Range("D36").Value = '2.658'
$ws.Range("E36").Value = '  -7.50%  '

$ws.Range("D37").Value = '0.9916'
$ws.Range("E37").Value = '  -9.59%  '

$ws.Range("D38").Value = '0.6036'
$ws.Range("E38").Value = '  -13.09%  '

$ws.Range("D39").Value = '2.695'
$ws.Range("E39").Value = '  -11.57%  '

$ws.Range("D40").Value = '0.01557'
$ws.Range("E40").Value = '  -4.37%  '

$ws.Range("D41").Value = '1.935'
$ws.Range("E41").Value = '  -10.38%  '

$ws.Range("E42").Value = '  -0.07%  '

$ws.Range("B43").Value = 'PaxosStandard'
$ws.Range("C43").Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
$ws.Range("D43").Value = '1.002'
$ws.Range("E43").Value = '  -0.05%  '

$ws.Range("B44").Value = 'Quant'
$ws.Range("C44").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D44").Value = '103.04'
$ws.Range("E44").Value = '  -1.52%  '

$ws.Range("B45").Value = 'TheSandbox'
$ws.Range("C45").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D45").Value = '0.3808'
$ws.Range("E45").Value = '  -14.00%  '

$ws.Range("B46").Value = 'TrustWalletToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D46").Value = '0.7405'
$ws.Range("E46").Value = '  -13.01%  '

$ws.Range("B47").Value = 'FraxShare'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D47").Value = '4.918'
$ws.Range("E47").Value = '  -13.65%  '

$ws.Range("B48").Value = 'Cronos'
$ws.Range("C48").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D48").Value = '0.05473'
$ws.Range("E48").Value = '  -1.35%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").Value = '0.1102'
$ws.Range("E49").Value = '  -5.46%  '

$ws.Range("B50").Value = 'Aptos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D50").Value = '5.956'
$ws.Range("E50").Value = '  -14.49%  '

$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '7.675'
$ws.Range("E51").Value = '  -10.10%  '
